$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.984.15"
$ws.Range("E2").Value = "  +1.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.688.88"
$ws.Range("E3").Value = "  +1.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.79"
$ws.Range("E5").Value = "  +0.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.80"
$ws.Range("E6").Value = "  +0.73%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -0.29%  "

# Row 9
$ws.Range("E9").Value = "  +2.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.98"
$ws.Range("E10").Value = "  +3.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.399"
$ws.Range("E11").Value = "  -3.74%  "

# Row 12
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000206"
$ws.Range("E13").Value = "  +7.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.90"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.171.23"
$ws.Range("E15").Value = "  +1.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.779.72"
$ws.Range("E16").Value = "  +1.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.681.12"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.70"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.86"
$ws.Range("E19").Value = "  -1.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  +3.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.17"
$ws.Range("E21").Value = "  -1.43%  "

# Row 22
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.82"
$ws.Range("E23").Value = "  +2.13%  "

# Row 24
$ws.Range("E24").Value = "  +15.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.91"
$ws.Range("E25").Value = "  +4.61%  "

# Row 26
$ws.Range("E26").Value = "  -5.70%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.172"
$ws.Range("E27").Value = "  +3.67%  "

# Row 28
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.65"
$ws.Range("E28").Value = "  -0.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  -1.59%  "

# Row 31
$ws.Range("E31").Value = "  +0.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.20"
$ws.Range("E32").Value = "  -5.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("E33").Value = "  -2.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.60"
$ws.Range("E34").Value = "  +2.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  -3.66%  "

# Row 36
$ws.Range("E36").Value = "  -0.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.65"
$ws.Range("E37").Value = "  +0.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.64"
$ws.Range("E38").Value = "  -0.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.99"
$ws.Range("E39").Value = "  -1.63%  "

# Row 40
$ws.Range("E40").Value = "  +0.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.44"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.84"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.13"
$ws.Range("E44").Value = "  -1.65%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0630"
$ws.Range("E45").Value = "  +0.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.53"
$ws.Range("E46").Value = "  +1.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  -1.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0265"
$ws.Range("E48").Value = "  +0.41%  "

# Row 49
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.68"
$ws.Range("E50").Value = "  +4.89%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0993"
$ws.Range("E51").Value = "  +0.94%  "
